{"js": "// Replace the literal placeholder-delimiter runs \"{{\" and \"}}\" with the\n// single-character brackets \"[\" and \"]\" respectively, everywhere they\n// occur in the document body (they appear as their own runs, e.g.\n// \"{{\" + \"\u7ba1\u7406\u756a\u53f7\" + \"}}\").\nconst openHits = context.document.body.search(\"{{\", { matchCase: true });\nopenHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of openHits.items) {\n  hit.insertText(\"[\", \"Replace\");\n}\nawait context.sync();\n\nconst closeHits = context.document.body.search(\"}}\", { matchCase: true });\ncloseHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of closeHits.items) {\n  hit.insertText(\"]\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace the literal placeholder-delimiter runs \"{{\" and \"}}\" with the\n# single-character brackets \"[\" and \"]\" respectively, everywhere they\n# occur in the document (they appear as their own runs, e.g.\n# \"{{\" + \"\u7ba1\u7406\u756a\u53f7\" + \"}}\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"{{\"\n$find.Replacement.Text = \"[\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"}}\"\n$find2.Replacement.Text = \"]\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2) | Out-Null\n"}
